$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from H1 (bold/centered/bordered) onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Data values for I2:I68 and J2:J68 (rows 2-68 correspond to data rows 2-68)
$iValues = @(7,8,7,9,8,7,7,8,7,6,8,7,6,7,6,6,7,8,10,7,7,5,7,7,8,8,11,8,7,6,6,8,7,7,7,9,5,9,7,5,6,8,6,7,7,6,7,7,8,7,8,7,7,7,7,7,7,7,6,4,7,5,6,7,5,8,4)
$jValues = @(7,8,8,9,9,7,7,8,7,6,8,7,6,7,6,6,7,8,10,7,7,6,7,7,8,8,11,8,7,6,6,8,7,7,7,9,5,9,7,6,6,8,6,7,7,6,7,7,8,7,8,7,7,7,7,7,7,7,7,5,7,5,6,7,5,8,4)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
